$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pour la prochaine fois")

# Remove the rows that disappeared from the "Pour la prochaine fois" sheet.
# Old layout:
#   B3=" "                                     B4="back/utils/socket..."   B5="Faire en sorte..."
#   B8="POV d'un joueur bug desfois..."
#   B10="optimiser le code en general"         C10="pour une prochiane co..."
#   B12="Sauvegarder kd par joueurs ? "
#   B14="Faille xss/autre securités"
#   B16="faire en sorte que le site..."
#   B18="quand on tue quelqu'un..."            C18="fait pour le joueur..."
#
# New layout drops B4, B5, C10, B12 and adds a new row 13 with a fresh task.
$ws.Range("B4").Value = $null
$ws.Range("B5").Value = $null
$ws.Range("C10").Value = $null
$ws.Range("B12").Value = $null

$ws.Range("B13").Value = "Refaire page d'accueil"

# Update the selection to match the saved cursor position.
$ws.Range("D8").Select()
